$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (shifts existing row 2 -> 3, row 3 -> 4)
$ws.Rows.Item(2).Insert()

# The inserted row inherits formatting from the row above (header row) by
# default; clear that so the new data row matches the plain, unstyled data
# rows below it.
$ws.Rows.Item(2).ClearFormats()

# Populate the new row 2 with the CRUDEOIL entry
$ws.Cells.Item(2, 1).Value = 435823
$ws.Cells.Item(2, 2).Value = "CRUDEOIL"
